$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current "Total Expenses" row (row 28) to hold
# the new Amazon Order 12 expense line. This shifts the old rows 28-30 down
# to 29-31.
$ws.Rows("28:28").Insert()

# Populate the new expense row (row 28)
$ws.Range("A28").Value = "Amazon Order 12"
$ws.Range("B28").Value = 43158
$ws.Range("C28").Value = "Morgan"
$ws.Range("D28").Value = "Amazon Order 12.pdf"
$ws.Range("E28").Value = 27.97
$ws.Range("F28").Value = "12V Power Adapters"

# Hyperlink the receipt description to the order PDF, matching the style
# used by the other receipt links in column D.
$ws.Hyperlinks.Add($ws.Range("D28"), "Amazon Order 12.pdf")
$ws.Range("D28").Style = "Hyperlink"

# Extend the running total to include the newly inserted row.
$ws.Range("E29").Formula = "=SUM(E2:E28)"

# Match the saved selection state from the source edit.
$ws.Range("A29").Select()
